$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = 86
$ws.Range("G3").Value = 14

# Row 4
$ws.Range("F4").Value = 63
$ws.Range("G4").Value = 18

# Row 14
$ws.Range("F14").Value = 80
$ws.Range("G14").Value = 20

# Row 17
$ws.Range("F17").Value = 73

# Row 20
$ws.Range("F20").Value = 71
$ws.Range("G20").Value = 29

# Row 21
$ws.Range("F21").Value = 84

# Row 22
$ws.Range("F22").Value = 85

# Row 24
$ws.Range("F24").Value = 49

# Row 25
$ws.Range("F25").Value = 92

# Row 27
$ws.Range("F27").Value = 77

# Row 28
$ws.Range("F28").Value = 89

# Row 30
$ws.Range("F30").Value = 97
$ws.Range("G30").Value = 3

# Row 31
$ws.Range("F31").Value = 81
$ws.Range("G31").Value = 19

# Row 32
$ws.Range("F32").Value = 73
$ws.Range("G32").Value = 27

# Row 33
$ws.Range("F33").Value = 74

# Row 34
$ws.Range("F34").Value = 84
$ws.Range("G34").Value = 16

# Row 35
$ws.Range("F35").Value = 72

# Row 37
$ws.Range("F37").Value = 87
$ws.Range("G37").Value = 13

# Row 38
$ws.Range("F38").Value = 78
